$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 397.57144
$ws.Range("I5").Value = 330.16666
$ws.Range("J5").Value = 802
$ws.Range("K5").Value = 330.16666
$ws.Range("L5").Value = 802
$ws.Range("M5").Value = -215.16666
$ws.Range("N5").Value = -1032

$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()

$ws.Range("H23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()

$ws.Range("H132").Value = 20768.83
$ws.Range("I132").Value = 21554.824
$ws.Range("J132").Value = 726
$ws.Range("K132").Value = 64664.472
$ws.Range("L132").Value = 2178
$ws.Range("M132").Value = -62134.472
$ws.Range("N132").Value = -7238

$ws.Range("H133").Value = 49225.22
$ws.Range("J133").Value = 49225.22
$ws.Range("L133").Value = 49225.22
$ws.Range("N133").Value = -59345.22

$ws.Range("H137").Value = 41669188
$ws.Range("I137").Value = 62501390
$ws.Range("K137").Value = 187504170
$ws.Range("M137").Value = -187501620

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 1117.6364
$ws.Range("I4").Value = 1179.4
$ws.Range("J4").Value = 500
$ws.Range("K4").Value = 1179.4
$ws.Range("L4").Value = 500
$ws.Range("M4").Value = -1063.4
$ws.Range("N4").Value = -732

$ws.Range("H19").Value = 4429.6
$ws.Range("I19").Value = 4132.6665
$ws.Range("J19").Value = 4875
$ws.Range("K19").Value = 4132.6665
$ws.Range("L19").Value = 4875
$ws.Range("M19").Value = -3903.6665
$ws.Range("N19").Value = -5333

$ws.Range("H32").Value = 4246.45
$ws.Range("I32").Value = 4246.45
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 4246.45
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -3959.45
$ws.Range("N32").ClearContents()

$ws.Range("H33").Value = 13000
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").ClearContents()

$ws.Range("H45").Value = 1168.8
$ws.Range("I45").Value = 1159.7778
$ws.Range("K45").Value = 1159.7778
$ws.Range("M45").Value = -782.7778000000001

$ws.Range("H55").Value = 0
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("M55").ClearContents()
$ws.Range("N55").ClearContents()

$ws.Range("H80").Value = 100000
$ws.Range("I80").Value = 100000
$ws.Range("K80").Value = 100000
$ws.Range("M80").Value = -99002

$ws.Range("H83").Value = 100000
$ws.Range("I83").Value = 100000
$ws.Range("K83").Value = 300000
$ws.Range("M83").Value = -295008

$ws.Range("H97").Value = 4477.68
$ws.Range("I97").Value = 4780.087
$ws.Range("K97").Value = 4780.087
$ws.Range("M97").Value = -4284.087

$ws.Range("H110").Value = 111112760
$ws.Range("I110").Value = 142858580
$ws.Range("J110").Value = 2406.5
$ws.Range("K110").Value = 142858580
$ws.Range("L110").Value = 2406.5
$ws.Range("M110").Value = -142856535
$ws.Range("N110").Value = -6496.5

$ws.Range("H122").Value = 2746.8
$ws.Range("I122").Value = 2247.125
$ws.Range("J122").Value = 3317.8572
$ws.Range("K122").Value = 6741.375
$ws.Range("L122").Value = 9953.571599999999
$ws.Range("M122").Value = -4291.375
$ws.Range("N122").Value = -14853.5716

$ws.Range("H133").Value = 47000
$ws.Range("J133").Value = 47000
$ws.Range("L133").Value = 47000
$ws.Range("N133").Value = -52060

$ws.Range("H139").Value = 48828.8
$ws.Range("J139").Value = 48828.8
$ws.Range("L139").Value = 48828.8
$ws.Range("N139").Value = -59108.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H59").Value = 59000
$ws.Range("J59").Value = 59000
$ws.Range("L59").Value = 59000
$ws.Range("N59").Value = -60694

$ws.Range("H107").Value = 2566.4666
$ws.Range("I107").Value = 2566.4666
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 2566.4666
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -646.4666000000002
$ws.Range("N107").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 56
$ws.Range("I7").Value = 56.333332
$ws.Range("J7").Value = 55
$ws.Range("K7").Value = 56.333332
$ws.Range("L7").Value = 55
$ws.Range("M7").Value = 56.666668
$ws.Range("N7").Value = -281

$ws.Range("H99").Value = 33335888
$ws.Range("I99").Value = 2628.5
$ws.Range("J99").Value = 166668930
$ws.Range("K99").Value = 2628.5
$ws.Range("L99").Value = 166668930
$ws.Range("M99").Value = -1130.5
$ws.Range("N99").Value = -166671926

$ws.Range("H122").Value = 2322.818
$ws.Range("I122").Value = 1563.7
$ws.Range("J122").Value = 9914
$ws.Range("K122").Value = 4691.1
$ws.Range("L122").Value = 29742
$ws.Range("M122").Value = -2241.1
$ws.Range("N122").Value = -34642

$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

$ws.Range("H126").Value = 33335888
$ws.Range("I126").Value = 2628.5
$ws.Range("J126").Value = 166668930
$ws.Range("K126").Value = 7885.5
$ws.Range("L126").Value = 500006790
$ws.Range("M126").Value = -5415.5
$ws.Range("N126").Value = -500011730

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 117.888885
$ws.Range("I11").Value = 152.2
$ws.Range("J11").Value = 75
$ws.Range("K11").Value = 456.6
$ws.Range("L11").Value = 225
$ws.Range("M11").Value = -316.6
$ws.Range("N11").Value = -505

$ws.Range("H12").Value = 76.22727
$ws.Range("I12").Value = 100.7
$ws.Range("J12").Value = 55.833332
$ws.Range("K12").Value = 302.1
$ws.Range("L12").Value = 167.499996
$ws.Range("M12").Value = -129.1
$ws.Range("N12").Value = -513.499996

$ws.Range("H131").Value = 8335134.5
$ws.Range("I131").Value = 943.3333
$ws.Range("J131").Value = 9010880
$ws.Range("K131").Value = 2829.9999
$ws.Range("L131").Value = 27032640
$ws.Range("M131").Value = 2210.0001
$ws.Range("N131").Value = -27042720

$ws.Range("H132").Value = 1493.5834
$ws.Range("I132").Value = 1374
$ws.Range("J132").Value = 1553.375
$ws.Range("K132").Value = 12366
$ws.Range("L132").Value = 13980.375
$ws.Range("M132").Value = -9836
$ws.Range("N132").Value = -19040.375

$ws.Range("H136").Value = 2822.6333
$ws.Range("J136").Value = 2863
$ws.Range("L136").Value = 8589
$ws.Range("N136").Value = -18789

$ws.Range("H140").Value = 7104.2046
$ws.Range("I140").Value = 11655.526
$ws.Range("J140").Value = 3645.2
$ws.Range("K140").Value = 34966.578
$ws.Range("L140").Value = 10935.6
$ws.Range("M140").Value = -29786.578
$ws.Range("N140").Value = -21295.6

$ws.Range("H141").Value = 6597.143
$ws.Range("I141").Value = 8045
$ws.Range("J141").Value = 4666.6665
$ws.Range("K141").Value = 24135
$ws.Range("L141").Value = 13999.9995
$ws.Range("M141").Value = -18955
$ws.Range("N141").Value = -24359.9995

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 79.333336
$ws.Range("I2").Value = 97.666664
$ws.Range("J2").Value = 61
$ws.Range("K2").Value = 97.666664
$ws.Range("L2").Value = 61
$ws.Range("M2").Value = 15.333336
$ws.Range("N2").Value = -287

$ws.Range("H122").Value = 2927.8333
$ws.Range("I122").Value = 2969.3076
$ws.Range("J122").Value = 2820
$ws.Range("K122").Value = 8907.9228
$ws.Range("L122").Value = 8460
$ws.Range("M122").Value = -6457.9228
$ws.Range("N122").Value = -13360

$ws.Range("H138").Value = 86476.336
$ws.Range("J138").Value = 86476.336
$ws.Range("L138").Value = 86476.336
$ws.Range("N138").Value = -96756.336

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 556.1
$ws.Range("I46").Value = 487.2857
$ws.Range("J46").Value = 716.6667
$ws.Range("K46").Value = 487.2857
$ws.Range("L46").Value = 716.6667
$ws.Range("M46").Value = -299.2857
$ws.Range("N46").Value = -1092.6667

$ws.Range("H55").Value = 493.44446
$ws.Range("I55").Value = 462.7143
$ws.Range("J55").Value = 601
$ws.Range("K55").Value = 462.7143
$ws.Range("L55").Value = 601
$ws.Range("M55").Value = -289.7143
$ws.Range("N55").Value = -947

$ws.Range("H61").Value = 1130.4
$ws.Range("I61").Value = 829.1429000000001
$ws.Range("J61").Value = 1833.3334
$ws.Range("K61").Value = 829.1429000000001
$ws.Range("L61").Value = 1833.3334
$ws.Range("M61").Value = -627.1429000000001
$ws.Range("N61").Value = -2237.3334

$ws.Range("H80").Value = 29333
$ws.Range("J80").Value = 29333
$ws.Range("L80").Value = 29333
$ws.Range("N80").Value = -31579

$ws.Range("H83").Value = 29333
$ws.Range("J83").Value = 29333
$ws.Range("L83").Value = 87999
$ws.Range("N83").Value = -99231

$ws.Range("H113").Value = 1130.4
$ws.Range("I113").Value = 829.1429000000001
$ws.Range("J113").Value = 1833.3334
$ws.Range("K113").Value = 829.1429000000001
$ws.Range("L113").Value = 1833.3334
$ws.Range("M113").Value = 1340.8571
$ws.Range("N113").Value = -6173.3334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 799.1111
$ws.Range("I113").Value = 838.93335
$ws.Range("J113").Value = 600
$ws.Range("K113").Value = 2516.80005
$ws.Range("L113").Value = 1800
$ws.Range("M113").Value = -346.8000499999998
$ws.Range("N113").Value = -6140

$ws.Range("H132").Value = 2795.9016
$ws.Range("I132").Value = 2411.6223
$ws.Range("J132").Value = 3876.6875
$ws.Range("K132").Value = 7234.8669
$ws.Range("L132").Value = 11630.0625
$ws.Range("M132").Value = -4704.8669
$ws.Range("N132").Value = -16690.0625
